$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.265.84'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '3.562.85'
$ws.Range("E3").Value = '  +4.27%  '
$ws.Range("E4").Value = '  +0.10%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.01'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +1.86%  '
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.44'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  +2.64%  '
$ws.Range("D7").Value = '3.562.88'
$ws.Range("E7").Value = '  +4.21%  '
$ws.Range("E8").Value = '  +0.18%  '
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.484'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  +3.76%  '
$ws.Range("E10").Value = '  +1.11%  '
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '8.04'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  +2.55%  '
$ws.Range("E12").Value = '  +1.99%  '
$ws.Range("D13").Value = '4.168.30'
$ws.Range("E13").Value = '  +4.04%  '
$ws.Range("E14").Value = '  +3.14%  '
$ws.Range("E15").Value = '  +1.30%  '
$ws.Range("D16").Value = '3.563.19'
$ws.Range("E16").Value = '  +3.76%  '
$ws.Range("D17").Value = '66.394.14'
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("E18").Value = '  -0.58%  '
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.50'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  +9.49%  '
$ws.Range("E20").Value = '  +2.25%  '
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.98'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  +0.97%  '
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '432.02'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  +3.84%  '
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.608'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  +4.74%  '
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.52'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  +0.81%  '
$ws.Range("D25").Value = '3.709.26'
$ws.Range("E25").Value = '  +4.20%  '
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("E27").Value = '  +8.99%  '
$ws.Range("E28").Value = '  +3.82%  '
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.07'
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = '  +3.36%  '
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("E31").Value = '  -0.21%  '
$ws.Range("E32").Value = '  +2.12%  '
$ws.Range("E33").Value = '  -0.78%  '
$ws.Range("D34").Value = '3.560.91'
$ws.Range("E34").Value = '  +4.10%  '
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.42'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  +3.50%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("E37").Value = '  +2.70%  '
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.92'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  +3.85%  '
$ws.Range("E39").Value = '  +2.43%  '
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  +0.00%  '
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '171.59'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  -0.96%  '
$ws.Range("E42").Value = '  -0.69%  '
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.23'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  +3.15%  '
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.897'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  +3.26%  '
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.96'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  +2.30%  '
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.98'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  +0.82%  '
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.50'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  -1.88%  '
$ws.Range("E48").Value = '  +4.40%  '
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.40'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  +5.67%  '
$ws.Range("E50").Value = '  +1.10%  '
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.953'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  +3.03%  '
